# Applies the "Add files via upload" change:
#   - Populates the previously-empty Lavatory / Pax seat / Galley sheets
#     (sheet2 / sheet3 / sheet4) with the same A/C-PN-DESCRIPTION-NOTE grid
#     that already exists on the Overhead sheet (sheet1), but with the
#     Description / PN / Note columns left blank.
#   - Copies the Overhead sheet's column widths / cell formatting onto the
#     three new sheets.
#   - Updates the selections on every sheet and makes "Galley" the active tab.

$wb = $excel.ActiveWorkbook

$wsOverhead = $wb.Worksheets.Item("Overhead")
$wsLavatory = $wb.Worksheets.Item("Lavatory")
$wsPaxSeat  = $wb.Worksheets.Item("Pax seat")
$wsGalley   = $wb.Worksheets.Item("Galley")

# A/C values that fill column A for rows 2-17 (same pattern used on Overhead).
$acByRow = @{
    2 = "B787"; 3 = "B787"; 4 = "B787"; 5 = "B787"; 6 = "B787"; 7 = "B787";
    8 = "ATR72"; 9 = "ATR72"; 10 = "ATR72"; 11 = "ATR72";
    12 = "A321"; 13 = "A321"; 14 = "A321"; 15 = "A321"; 16 = "A321"; 17 = "A321"
}

foreach ($ws in @($wsLavatory, $wsPaxSeat, $wsGalley)) {
    # Bring over the Overhead formatting (fonts / borders / number formats)
    # for A1:D17 so the new sheets reuse the same cell styles instead of
    # minting new ones.
    $wsOverhead.Range("A1:D17").Copy()
    $ws.Range("A1:D17").PasteSpecial(-4122)  # xlPasteFormats

    # Header row.
    $ws.Range("A1").Value = "A/C"
    $ws.Range("B1").Value = "DESCRIPTION"
    $ws.Range("C1").Value = "PART NUMBER (PN)"
    $ws.Range("D1").Value = "NOTE"

    # Column A data rows; B/C/D stay blank (formatting only).
    foreach ($r in 2..17) {
        $ws.Cells.Item($r, 1).Value = $acByRow[$r]
    }

    # Column widths matching the Overhead sheet.
    $ws.Columns.Item(1).ColumnWidth = 14.140625
    $ws.Columns.Item(2).ColumnWidth = 20.28515625
    $ws.Columns.Item(3).ColumnWidth = 23
    $ws.Columns.Item(4).ColumnWidth = 46.42578125

    $ws.Range("B2:D3").Select()
}

# Overhead's own selection moves to the full used range.
$wsOverhead.Range("A1:D17").Select()

# Galley becomes the active/visible tab.
$wsGalley.Activate()
$wsGalley.Range("B2:D3").Select()
